# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# to the latest scraped values, per the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''24.880.10'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '''1.709.88'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''312.00'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '''0.9984'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '''0.3764'
$ws.Range("E7").Value = '  +1.31%  '
$ws.Range("D8").Value = '''49.61'
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("D9").Value = '''0.3453'
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").Value = '''1.213'
$ws.Range("E10").Value = '  +2.39%  '
$ws.Range("D11").Value = '''0.07557'
$ws.Range("E11").Value = '  +3.85%  '
$ws.Range("D12").Value = '''0.9993'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '''21.28'
$ws.Range("E13").Value = '  +4.11%  '
$ws.Range("D14").Value = '''6.333'
$ws.Range("E14").Value = '  +3.19%  '
$ws.Range("D15").Value = '''7.083'
$ws.Range("E15").Value = '  +4.94%  '
$ws.Range("D16").Value = '''1.710.19'
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").Value = '''0.00001136'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").Value = '''0.06729'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").Value = '''0.9986'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '''85.15'
$ws.Range("E20").Value = '  +4.67%  '
$ws.Range("D21").Value = '''17.38'
$ws.Range("E21").Value = '  +5.54%  '
$ws.Range("D22").Value = '''6.418'
$ws.Range("E22").Value = '  +5.08%  '
$ws.Range("D23").Value = '''13.33'
$ws.Range("E23").Value = '  +11.20%  '
$ws.Range("D24").Value = '''24.866.85'
$ws.Range("E24").Value = '  +1.94%  '
$ws.Range("D25").Value = '''2.455'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").Value = '''2.811'
$ws.Range("E26").Value = '  +5.13%  '
$ws.Range("E27").Value = '  +5.07%  '
$ws.Range("D28").Value = '''152.07'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = '''133.20'
$ws.Range("E29").Value = '  +5.01%  '
$ws.Range("D30").Value = '''1.900.76'
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").Value = '''1.253'
$ws.Range("E31").Value = '  +28.67%  '
$ws.Range("D32").Value = '''6.953'
$ws.Range("E32").Value = '  +9.38%  '
$ws.Range("D33").Value = '''4.250'
$ws.Range("E33").Value = '  +5.38%  '
$ws.Range("D34").Value = '''13.98'
$ws.Range("E34").Value = '  +12.32%  '
$ws.Range("D35").Value = '''1.795'
$ws.Range("E35").Value = '  +7.43%  '
$ws.Range("D36").Value = '''0.08817'
$ws.Range("E36").Value = '  +4.09%  '
$ws.Range("D37").Value = '''9.416'
$ws.Range("E37").Value = '  +5.19%  '
$ws.Range("D38").Value = '''5.635'
$ws.Range("E38").Value = '  +5.18%  '
$ws.Range("D39").Value = '''0.06705'
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("D40").Value = '''0.02424'
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("D41").Value = '''0.2244'
$ws.Range("E41").Value = '  +6.14%  '
$ws.Range("D42").Value = '''1.285'
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D43").Value = '''0.6469'
$ws.Range("E43").Value = '  +4.39%  '
$ws.Range("D44").Value = '''0.9982'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '''14.03'
$ws.Range("E45").Value = '  +8.09%  '
$ws.Range("D46").Value = '''0.6185'
$ws.Range("E46").Value = '  +3.80%  '
$ws.Range("D47").Value = '''3.829'
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("E48").Value = '  +5.12%  '
$ws.Range("D49").Value = '''130.54'
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("D50").Value = '''0.07324'
$ws.Range("D51").Value = '''80.30'
$ws.Range("E51").Value = '  +5.99%  '
